$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Body text replacements (Simplified -> Traditional Chinese + wording tweaks)
Replace-Text "英语" "英語"
Replace-Text " / 葡萄牙语 / 法语 / 泰语 / 越南语 / 西班牙语" " / 葡萄牙語 / 法語 / 泰語 / 越南語 / 西班牙語"
Replace-Text "简介" "簡介"
Replace-Text "一封发送给目标国家已确认出席的合作伙伴的电子邮件。 我们希望他们提交他们的文件。 将通过 customer.io 发送" "發送給目標國家已回覆參加的合作夥伴的電子郵件。 我們希望他們提交他們的文件。 將通過 customer.io 發送。"
Replace-Text "目标受众" "目標受眾"
Replace-Text "被邀请且已确认出席的合作伙伴" "已邀請並確認參加的合作夥伴"
Replace-Text "主题: " "主題:"
Replace-Text "[活动名称]" "[活動名稱]"
Replace-Text " — 下一步" " — 進一步行動"
Replace-Text "感谢您报名参加 " "感謝您報名參加 "
Replace-Text "[合作伙伴姓名]" "[合作夥伴姓名]"
Replace-Text "很高兴您能参加 " "很高興您能參加 "
Replace-Text "! " "！ "
Replace-Text "为了确认您的注册，需要您和您选择的一位嘉宾向我们提供：" "為了確認您的註冊，需要您和您選擇的一位嘉賓向我們提供："
Replace-Text "经签名的 " "經簽名的 "
Replace-Text "《行为准则》" "行為守則 "
Replace-Text "《条款和条件》" "條款和條件"
Replace-Text "（每人 1 份）" "（每人 1 份）("
Replace-Text "国际护照扫描件" "國際護照掃描件"
Replace-Text "Covid-19 疫苗接种证书" "Covid-19 疫苗接種證書"
Replace-Text "发送我的详细信息" "發送我的詳細資料"
Replace-Text "区域经理将与您联系以确认您的预订或索取任何其他相关详细信息。" "區域經理將與您聯繫以確認預訂或索取任何其他相關詳細資料。 "
Replace-Text "我们的活动套餐为您和您的嘉宾提供：" "我們的活動套餐為您和您的嘉賓提供： "
Replace-Text "机票 " "機票 "
Replace-Text "旅行保险 " "旅遊保險 "
Replace-Text "机场-酒店-机场接送 " "機場 — 酒店 — 機場接送 "
Replace-Text "为您和您的嘉宾提供一间酒店客房/为您和您的嘉宾提供两间酒店客房" "為您和您的嘉賓提供一間酒店客房/為您和您的嘉賓提供兩間酒店客房"
Replace-Text "观光游览 " "觀光遊覽 "
Replace-Text "在您的出发日期之前，我们将向您发送确认信，其中包含活动日程以及航班、交通和住宿信息。" "將在出發日期前向您傳送確認信，其中包含活動日程以及航班、交通和住宿資訊。 "
Replace-Text "如果您有任何疑问，请通过 " "如有任何疑問，請與我們聯繫："
Replace-Text "实时聊天" "即時聊天"
Replace-Text " 联系我们。" "。 "
Replace-Text "如有疑问，请联系我们的区域经理 " "如有任何疑問，請聯繫您的區域經理 "
Replace-Text "[名字]" "[姓名]"
Replace-Text "，邮箱 " "，可通過 "
Replace-Text "[电子邮件地址]" "[電子郵件地址]"
Replace-Text "[WHATSAPP 号码]" "[WHATSAPP 號碼]"
Replace-Text " (WhatsApp)。 " " (WhatsApp) 聯繫。 "
Replace-Text "期待很快与您见面。" "期待很快與您見面。"

# Comment text replacements (by comment id, via Comments collection index = id + 1)
$d.Comments.Item(3).Range.Text = "鏈接到條款和條件"
$d.Comments.Item(2).Range.Text = "鏈接到行為守則"
$d.Comments.Item(4).Range.Text = "請確認這些"
$d.Comments.Item(5).Range.Text = "選擇其中一個"
$d.Comments.Item(1).Range.Text = "請檢查這些是否為所有所需文件"
